$d = $word.ActiveDocument

# 1) Fix "tests results" -> "test results" typo in the week-3 log entry.
$d.Content.Find.Execute("tests results from the library session", $true, $false, $false, $false, $false, $true, 1, $false, "test results from the library session", 2) | Out-Null

# 2) Insert a new log entry paragraph before the "LEDR[2:0] ..." paragraph.
$find = $d.Content
$found = $find.Find.Execute("             - LEDR[2:0] have been used to signal the status of different inputs.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $insertRange = $d.Range($find.Start, $find.Start)
    $insertRange.InsertParagraphBefore()
    $newPara = $d.Range($find.Start, $find.Start)
    $newPara.Text = "             - The number of lives HEX display panel has been tested on the screen. All 3 panels are functioning. Some refreshing problems have been noticed."
    $newPara.ParagraphFormat.Style = "Normal1"
}
